# Update the three-digit-divided-by-one-digit division problems in the
# document's tables to the new set of values, per the commit diff.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "879÷7="; New = "479÷3=" },
    @{ Old = "575÷8="; New = "897÷4=" },
    @{ Old = "368÷8="; New = "532÷9=" },
    @{ Old = "582÷9="; New = "375÷2=" },
    @{ Old = "300÷2="; New = "253÷5=" },
    @{ Old = "755÷4="; New = "885÷6=" },
    @{ Old = "685÷9="; New = "889÷2=" },
    @{ Old = "670÷8="; New = "837÷5=" },
    @{ Old = "937÷4="; New = "783÷2=" },
    @{ Old = "837÷9="; New = "174÷2=" },
    @{ Old = "111÷3="; New = "390÷2=" },
    @{ Old = "538÷4="; New = "437÷6=" },
    @{ Old = "780÷2="; New = "306÷6=" },
    @{ Old = "430÷9="; New = "390÷8=" },
    @{ Old = "189÷7="; New = "585÷2=" },
    @{ Old = "529÷7="; New = "132÷4=" },
    @{ Old = "107÷7="; New = "217÷7=" },
    @{ Old = "152÷2="; New = "586÷7=" },
    @{ Old = "355÷6="; New = "384÷3=" },
    @{ Old = "802÷5="; New = "995÷8=" },
    @{ Old = "534÷9="; New = "958÷8=" },
    @{ Old = "811÷3="; New = "189÷9=" },
    @{ Old = "267÷6="; New = "341÷7=" },
    @{ Old = "373÷8="; New = "587÷3=" },
    @{ Old = "862÷3="; New = "790÷4=" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $pair.New, 2)
}

$d.Save()
